$wb = $excel.ActiveWorkbook

# Update the status text from "Ready for handoff" to "In Translation"
# across all sheets that reference it.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns("E:F").AutoFit() | Out-Null

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns("C:C").AutoFit() | Out-Null

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns("C:C").AutoFit() | Out-Null
